# Apply edits to menarche_2018 worksheet:
#  - Select cell E5 (update sheet view selection)
#  - Flip several "menarche" (column C) values between 0 and 1

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the current selection/active cell shown in the sheet view
$ws.Range("E5").Select()

# Rows where column C changes from 0 -> 1
$rowsToOne = @(6, 40, 41, 42, 43, 44, 70, 71)
foreach ($r in $rowsToOne) {
    $ws.Cells.Item($r, 3).Value = 1
}

# Rows where column C changes from 1 -> 0
$rowsToZero = @(55, 56)
foreach ($r in $rowsToZero) {
    $ws.Cells.Item($r, 3).Value = 0
}
